# Se mejora el formato fecha yyyy-mm-dd
# Actual content change captured by the diff: the header text in N1
# ("Población_X_1000") loses its accent, becoming "Poblacion_X_1000".
# (Excel drops the now-unused old shared string and appends the new
# text at the end of the shared-strings table, which is why every
# other shared-string index from 13 onward shifts down by one - that
# is an automatic side effect, not something to script by hand.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "Poblacion_X_1000"

# The active selection moves from C12 to E1.
$ws.Range("E1").Select()
